$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two weekly date blocks (rows 2-5 vs rows 6-9) have swapped places:
# what used to be the 44890 (La Ligua) block is now the 44908
# (Provincia de Limarí) block, and vice versa. Swap the per-row values
# for the columns that carry the date-specific data: D, M, N, O, P, R, S.

$pairs = @(
    @{ a = 2; b = 6 },
    @{ a = 3; b = 7 },
    @{ a = 4; b = 8 },
    @{ a = 5; b = 9 }
)

$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($pair in $pairs) {
    $rowA = $pair.a
    $rowB = $pair.b

    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}
